# Updated ITA model - 2025-09-01 12:44
$wb = $excel.ActiveWorkbook

# --- Update ev_charging_uc sheet: TimeSlices table now maps day_night/timeslice to AllSaAllH ---
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Cells.Item(12, 2).Value = "day_night"
$wsEv.Cells.Item(12, 3).Value = "timeslice"
$wsEv.Cells.Item(13, 2).Value = "D"
$wsEv.Cells.Item(13, 3).Value = "AllSaAllH"
$wsEv.Cells.Item(14, 2).Value = "N"
$wsEv.Cells.Item(14, 3).Value = "AllSaAllH"

# --- Add new sheet: re_profiles (after timeslice_def) ---
$wsTsDef = $wb.Worksheets.Item("timeslice_def")
$wsRe = $wb.Worksheets.Add($null, $wsTsDef)
$wsRe.Name = "re_profiles"

$wsRe.Cells.Item(9, 2).Value = "~TFM_DINS-AT"
$wsRe.Cells.Item(10, 2).Value = "commodity"
$wsRe.Cells.Item(10, 3).Value = "timeslice"
$wsRe.Cells.Item(10, 4).Value = "com_fr"
$wsRe.Cells.Item(10, 5).Value = "process"
$wsRe.Cells.Item(11, 2).Value = "elc_spv-ITA"
$wsRe.Cells.Item(11, 3).Value = "AllSaAllH"
$wsRe.Cells.Item(11, 4).Value = 0.99999999999980893
$wsRe.Cells.Item(11, 5).Value = "IMPNRGZ"

$wsRe.Cells.Item(9, 7).Value = "~TFM_DINS-AT"
$wsRe.Cells.Item(10, 7).Value = "commodity"
$wsRe.Cells.Item(10, 8).Value = "timeslice"
$wsRe.Cells.Item(10, 9).Value = "com_fr"
$wsRe.Cells.Item(10, 10).Value = "process"
$wsRe.Cells.Item(11, 7).Value = "elc_won-ITA"
$wsRe.Cells.Item(11, 8).Value = "AllSaAllH"
$wsRe.Cells.Item(11, 9).Value = 0.99999999999979439
$wsRe.Cells.Item(11, 10).Value = "IMPNRGZ"

$wsRe.Cells.Item(9, 12).Value = "~TFM_DINS-AT"
$wsRe.Cells.Item(10, 12).Value = "commodity"
$wsRe.Cells.Item(10, 13).Value = "timeslice"
$wsRe.Cells.Item(10, 14).Value = "com_fr"
$wsRe.Cells.Item(10, 15).Value = "process"
$wsRe.Cells.Item(11, 12).Value = "elc_wof-ITA"
$wsRe.Cells.Item(11, 13).Value = "AllSaAllH"
$wsRe.Cells.Item(11, 14).Value = 0
$wsRe.Cells.Item(11, 15).Value = "IMPNRGZ"

$wsRe.Cells.Item(9, 17).Value = "~TFM_INS-AT"
$wsRe.Cells.Item(10, 17).Value = "timeslice"
$wsRe.Cells.Item(10, 18).Value = "ncap_afs"
$wsRe.Cells.Item(10, 19).Value = "pset_ci"
$wsRe.Cells.Item(11, 17).Value = "AllS"
$wsRe.Cells.Item(11, 18).Value = 1.2
$wsRe.Cells.Item(11, 19).Value = "hydro"

# --- Add new sheet: load_shapes (after re_profiles) ---
$wsLoad = $wb.Worksheets.Add($null, $wsRe)
$wsLoad.Name = "load_shapes"

$wsLoad.Cells.Item(9, 2).Value = "~TFM_DINS-AT"
$wsLoad.Cells.Item(10, 2).Value = "g_yrfr"
$wsLoad.Cells.Item(10, 3).Value = "com_fr"
$wsLoad.Cells.Item(10, 4).Value = "timeslice"
$wsLoad.Cells.Item(10, 5).Value = "commodity"
$wsLoad.Cells.Item(11, 2).Value = 1
$wsLoad.Cells.Item(11, 3).Value = 1.0000000000000002
$wsLoad.Cells.Item(11, 4).Value = "AllSaAllH"
$wsLoad.Cells.Item(11, 5).Value = "elc_roadtransport"

$wsLoad.Cells.Item(9, 7).Value = "~TFM_DINS-AT"
$wsLoad.Cells.Item(10, 7).Value = "commodity"
$wsLoad.Cells.Item(10, 8).Value = "timeslice"
$wsLoad.Cells.Item(10, 9).Value = "com_fr"
$wsLoad.Cells.Item(11, 7).Value = "elc_buildings"
$wsLoad.Cells.Item(11, 8).Value = "AllSaAllH"
$wsLoad.Cells.Item(11, 9).Value = 1
$wsLoad.Cells.Item(12, 7).Value = "elc_industry"
$wsLoad.Cells.Item(12, 8).Value = "AllSaAllH"
$wsLoad.Cells.Item(12, 9).Value = 1

$wsLoad.Cells.Item(9, 11).Value = "~TFM_DINS-AT"
$wsLoad.Cells.Item(10, 11).Value = "commodity"
$wsLoad.Cells.Item(10, 12).Value = "timeslice"
$wsLoad.Cells.Item(10, 13).Value = "com_pkflx"
$wsLoad.Cells.Item(11, 11).Value = "ELC"
$wsLoad.Cells.Item(11, 12).Value = "AllSaAllH"
$wsLoad.Cells.Item(11, 13).Value = 0.29960038929139898

# --- Make load_shapes the active/selected sheet (last tab selected) ---
$wsLoad.Select()
